# Grading update, new warm up (Homework 5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: "Homework 5" header ---
$ws.Range("I1").Value = "Homework 5"
$ws.Columns.Item(9).ColumnWidth = 10.1

# --- New column I grades (score/total style formulas, like the other homework columns) ---
$ws.Range("I2").Formula  = "=25/25"
$ws.Range("I3").Formula  = "=25/25"
$ws.Range("I4").Formula  = "=25/25"
$ws.Range("I5").Formula  = "=0"
$ws.Range("I6").Formula  = "=0"
$ws.Range("I7").Formula  = "=20/25"
$ws.Range("I8").Formula  = "=25/25"
$ws.Range("I9").Formula  = "=20/25"
$ws.Range("I10").Formula = "=0"
$ws.Range("I11").Formula = "=0"
$ws.Range("I12").Formula = "=0"
$ws.Range("I13").Formula = "=25/25"
$ws.Range("I14").Formula = "=25/25"
$ws.Range("I15").Formula = "=25/25"

# --- Grading corrections on existing homework columns ---
# Nahom Anteneh (row 5): Homework 1 grade recorded
$ws.Range("D5").Formula = "=35/35"

# Kai Stephens (row 7): Homework 3 grade recorded
$ws.Range("F7").Formula = "=29/30"

# Almas Waseem (row 11): Homework 1-3 grades recorded
$ws.Range("D11").Formula = "=30/35"
$ws.Range("E11").Formula = "=21/25"
$ws.Range("F11").Formula = "=22/30"

# Ty Carlson (row 12): Homework 4 re-saved as an explicit formula (value unchanged)
$ws.Range("H12").Formula = "=0"

# --- Put the active selection on F8, matching the author's cursor position ---
$ws.Range("F8").Select()
